$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text (coin name / URL) - safe to set directly.
# Columns D and E contain numeric-looking / percentage text that Excel would
# otherwise auto-convert to numbers, so force General/"@" text formatting,
# write the literal text, then restore the default "Normal" style so the
# cell keeps the workbook's original (unstyled) appearance.

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.390.79"
Set-TextValue $ws.Range("E2") "  +0.00%  "
Set-TextValue $ws.Range("D3") "1.848.26"
Set-TextValue $ws.Range("E3") "  +0.05%  "
Set-TextValue $ws.Range("D4") "0.9995"
Set-TextValue $ws.Range("E4") "  +0.08%  "
Set-TextValue $ws.Range("D5") "240.27"
Set-TextValue $ws.Range("E5") "  -0.06%  "
Set-TextValue $ws.Range("D6") "0.6292"
Set-TextValue $ws.Range("E6") "  -0.13%  "
Set-TextValue $ws.Range("E7") "  +0.04%  "
Set-TextValue $ws.Range("D8") "0.07626"
Set-TextValue $ws.Range("E8") "  +1.17%  "
Set-TextValue $ws.Range("E9") "  -0.52%  "
Set-TextValue $ws.Range("D10") "24.48"
Set-TextValue $ws.Range("E10") "  -0.06%  "
Set-TextValue $ws.Range("D11") "0.07746"
Set-TextValue $ws.Range("E11") "  +0.32%  "
Set-TextValue $ws.Range("D12") "1.842.97"
Set-TextValue $ws.Range("E12") "  -0.38%  "
Set-TextValue $ws.Range("E13") "  +0.31%  "
Set-TextValue $ws.Range("D14") "0.00001090"
Set-TextValue $ws.Range("E14") "  +8.90%  "
Set-TextValue $ws.Range("D15") "0.6790"
Set-TextValue $ws.Range("E15") "  -0.63%  "
Set-TextValue $ws.Range("E16") "  +0.61%  "
Set-TextValue $ws.Range("D17") "2.092.48"
Set-TextValue $ws.Range("E17") "  -7.59%  "
Set-TextValue $ws.Range("D18") "6.131"
Set-TextValue $ws.Range("D19") "29.422.28"
Set-TextValue $ws.Range("E19") "  -0.03%  "
Set-TextValue $ws.Range("D20") "228.39"
Set-TextValue $ws.Range("E21") "  +0.11%  "
Set-TextValue $ws.Range("D22") "1.000"
Set-TextValue $ws.Range("D23") "7.443"
Set-TextValue $ws.Range("E23") "  -1.31%  "
Set-TextValue $ws.Range("E24") "  -0.02%  "
Set-TextValue $ws.Range("D25") "157.21"
Set-TextValue $ws.Range("E25") "  +0.10%  "
Set-TextValue $ws.Range("E26") "  -0.65%  "
Set-TextValue $ws.Range("D27") "8.372"
Set-TextValue $ws.Range("E27") "  -0.02%  "
Set-TextValue $ws.Range("E28") "  +0.01%  "
Set-TextValue $ws.Range("D29") "1.467"
Set-TextValue $ws.Range("E29") "  -0.04%  "
Set-TextValue $ws.Range("E30") "  +3.55%  "
Set-TextValue $ws.Range("D31") "0.05630"
Set-TextValue $ws.Range("E31") "  -0.97%  "
Set-TextValue $ws.Range("D32") "4.111"
Set-TextValue $ws.Range("E32") "  -0.49%  "
Set-TextValue $ws.Range("D33") "4.042"
Set-TextValue $ws.Range("E33") "  +0.57%  "
Set-TextValue $ws.Range("E34") "  +0.47%  "
Set-TextValue $ws.Range("D35") "1.156"
Set-TextValue $ws.Range("E35") "  +0.06%  "
Set-TextValue $ws.Range("D36") "0.7092"
Set-TextValue $ws.Range("E36") "  -0.69%  "
Set-TextValue $ws.Range("D37") "2.587"
Set-TextValue $ws.Range("E37") "  -0.04%  "
Set-TextValue $ws.Range("D38") "2.773"
Set-TextValue $ws.Range("E38") "  -0.34%  "
Set-TextValue $ws.Range("D39") "1.229.55"
Set-TextValue $ws.Range("E39") "  -1.94%  "
Set-TextValue $ws.Range("D40") "0.01795"
Set-TextValue $ws.Range("E40") "  -1.02%  "
Set-TextValue $ws.Range("D41") "6.485"
Set-TextValue $ws.Range("E41") "  +4.45%  "
Set-TextValue $ws.Range("D42") "0.9074"
Set-TextValue $ws.Range("E42") "  -0.60%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D44") "2.001.76"
Set-TextValue $ws.Range("E44") "  -0.01%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "101.46"
Set-TextValue $ws.Range("E45") "  +0.59%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "66.03"
Set-TextValue $ws.Range("E46") "  -0.30%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D47") "0.00000000122"
Set-TextValue $ws.Range("E47") "  +4.14%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D48") "7.149"
Set-TextValue $ws.Range("E48") "  +1.37%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D49") "0.4009"
Set-TextValue $ws.Range("E49") "  -0.35%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "9.009"
Set-TextValue $ws.Range("E50") "  -0.92%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D51") "1.684"
Set-TextValue $ws.Range("E51") "  -0.30%  "
